$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents()
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("H31").Value = 300
$ws.Range("I31").Value = 300
$ws.Range("K31").Value = 900
$ws.Range("M31").Value = -670

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5160.99
$ws.Range("I32").Value = 3534.2532
$ws.Range("J32").Value = 11280.619
$ws.Range("K32").Value = 3534.2532
$ws.Range("L32").Value = 11280.619
$ws.Range("M32").Value = -3247.2532
$ws.Range("N32").Value = -11854.619
$ws.Range("H44").Value = 29498.111
$ws.Range("J44").Value = 29498.111
$ws.Range("L44").Value = 29498.111
$ws.Range("N44").Value = -30474.111
$ws.Range("H55").Value = 34416.8
$ws.Range("J55").Value = 34416.8
$ws.Range("L55").Value = 34416.8
$ws.Range("N55").Value = -35046.8
$ws.Range("H109").Value = 25906.889
$ws.Range("J109").Value = 25906.889
$ws.Range("L109").Value = 25906.889
$ws.Range("N109").Value = -28680.889

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7096169
$ws.Range("I31").Value = 3926.8206
$ws.Range("K31").Value = 3926.8206
$ws.Range("M31").Value = -3631.8206
$ws.Range("H34").Value = 7096169
$ws.Range("I34").Value = 3926.8206
$ws.Range("K34").Value = 3926.8206
$ws.Range("M34").Value = -3724.8206
$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("M44").ClearContents()
$ws.Range("H68").Value = 19809.4
$ws.Range("J68").Value = 19809.4
$ws.Range("L68").Value = 19809.4
$ws.Range("N68").Value = -21307.4
$ws.Range("H71").Value = 19809.4
$ws.Range("J71").Value = 19809.4
$ws.Range("L71").Value = 59428.2
$ws.Range("N71").Value = -66916.20000000001
$ws.Range("H109").Value = 27530
$ws.Range("J109").Value = 27530
$ws.Range("L109").Value = 27530
$ws.Range("N109").Value = -29610
$ws.Range("H132").Value = 13515829
$ws.Range("I132").Value = 20001970
$ws.Range("K132").Value = 60005910
$ws.Range("M132").Value = -60003380

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H45").Value = 500
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 500
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 1500
$ws.Range("M45").ClearContents()
$ws.Range("N45").Value = -2564
$ws.Range("H61").Value = 290.7
$ws.Range("I61").Value = 91.4
$ws.Range("J61").Value = 490
$ws.Range("K61").Value = 274.2
$ws.Range("L61").Value = 1470
$ws.Range("M61").Value = -59.20000000000005
$ws.Range("N61").Value = -1900
$ws.Range("H98").Value = 934.5294
$ws.Range("I98").Value = 834.6
$ws.Range("J98").Value = 976.1667
$ws.Range("K98").Value = 2503.8
$ws.Range("L98").Value = 2928.5001
$ws.Range("M98").Value = -1005.8
$ws.Range("N98").Value = -5924.5001
$ws.Range("H105").Value = 6500
$ws.Range("J105").Value = 6500
$ws.Range("L105").Value = 19500
$ws.Range("N105").Value = -24742
$ws.Range("H109").Value = 4234.5
$ws.Range("I109").Value = 2014
$ws.Range("J109").Value = 4436.364
$ws.Range("K109").Value = 6042
$ws.Range("L109").Value = 13309.092
$ws.Range("M109").Value = -5002
$ws.Range("N109").Value = -15389.092
$ws.Range("H114").Value = 2570.9333
$ws.Range("I114").Value = 2594.4
$ws.Range("J114").Value = 2559.2
$ws.Range("K114").Value = 7783.200000000001
$ws.Range("L114").Value = 7677.599999999999
$ws.Range("M114").Value = -4529.200000000001
$ws.Range("N114").Value = -14185.6
$ws.Range("H116").Value = 618.8
$ws.Range("I116").Value = 523.5
$ws.Range("J116").Value = 1000
$ws.Range("K116").Value = 1570.5
$ws.Range("L116").Value = 3000
$ws.Range("M116").Value = 1871.5
$ws.Range("N116").Value = -9884
$ws.Range("H117").Value = 427
$ws.Range("I117").Value = 250
$ws.Range("J117").Value = 515.5
$ws.Range("K117").Value = 750
$ws.Range("L117").Value = 1546.5
$ws.Range("M117").Value = 2692
$ws.Range("N117").Value = -8430.5
$ws.Range("H118").Value = 1325.6666
$ws.Range("J118").Value = 1467.9
$ws.Range("L118").Value = 4403.700000000001
$ws.Range("N118").Value = -6889.700000000001
$ws.Range("H119").Value = 8997.833
$ws.Range("I119").Value = 5995.6665
$ws.Range("J119").Value = 12000
$ws.Range("K119").Value = 17986.9995
$ws.Range("L119").Value = 36000
$ws.Range("M119").Value = -13148.9995
$ws.Range("N119").Value = -45676

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 9271.286
$ws.Range("I57").Value = 3979.8
$ws.Range("J57").Value = 22500
$ws.Range("K57").Value = 3979.8
$ws.Range("L57").Value = 22500
$ws.Range("M57").Value = -3159.8
$ws.Range("N57").Value = -24140
$ws.Range("H118").Value = 14504.348
$ws.Range("J118").Value = 14504.348
$ws.Range("L118").Value = 14504.348
$ws.Range("N118").Value = -17818.348
$ws.Range("H132").Value = 3416.1833
$ws.Range("I132").Value = 3847.5227
$ws.Range("J132").Value = 2230
$ws.Range("K132").Value = 11542.5681
$ws.Range("L132").Value = 6690
$ws.Range("M132").Value = -9012.5681
$ws.Range("N132").Value = -11750

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 726.04
$ws.Range("I93").Value = 621.5
$ws.Range("J93").Value = 911.8889
$ws.Range("K93").Value = 621.5
$ws.Range("L93").Value = 911.8889
$ws.Range("M93").Value = 626.5
$ws.Range("N93").Value = -3407.8889
$ws.Range("H100").Value = 3259.625
$ws.Range("I100").Value = 3031.7778
$ws.Range("J100").Value = 3552.5715
$ws.Range("K100").Value = 3031.7778
$ws.Range("L100").Value = 3552.5715
$ws.Range("M100").Value = -2490.7778
$ws.Range("N100").Value = -4634.5715
$ws.Range("H136").Value = 13162356
$ws.Range("I136").Value = 13890542
$ws.Range("J136").Value = 55002.5
$ws.Range("K136").Value = 41671626
$ws.Range("L136").Value = 165007.5
$ws.Range("M136").Value = -41669076
$ws.Range("N136").Value = -170107.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 830.46875
$ws.Range("I136").Value = 822.6667
$ws.Range("J136").Value = 947.5
$ws.Range("K136").Value = 2468.0001
$ws.Range("L136").Value = 2842.5
$ws.Range("M136").Value = 81.9998999999998
$ws.Range("N136").Value = -7942.5
